$wb = $excel.ActiveWorkbook

$wsScrewed = $wb.Worksheets.Item("screwed")
$wsShifts  = $wb.Worksheets.Item("shiftsPerWorker")

# Append four new rows (12-15) to the "screwed" sheet
$wsScrewed.Range("A12").Value = 11
$wsScrewed.Range("B12").Value = "stav"

$wsScrewed.Range("A13").Value = 12
$wsScrewed.Range("B13").Value = "stav"

$wsScrewed.Range("A14").Value = 13
$wsScrewed.Range("B14").Value = "adir"

$wsScrewed.Range("A15").Value = 14
$wsScrewed.Range("B15").Value = "adir"

# Update the shift counts on the "shiftsPerWorker" sheet
$wsShifts.Range("B2").Value = 1
$wsShifts.Range("B3").Value = 6
$wsShifts.Range("B4").Value = 5
$wsShifts.Range("B5").Value = 3
$wsShifts.Range("B7").Value = 4
